$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# Column A holds a date-like label ("01-08-2021") that must be stored as
# plain text (shared string), not auto-converted into a date serial number.
# Assigning it through .Formula as a quoted string literal keeps it text,
# then Copy / PasteSpecial(values) bakes it back down to a literal text
# cell without leaving any NumberFormat/style residue behind.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Formula = '="01-08-2021"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

$values = @(
    8258.799999999999,
    471.4,
    215.6,
    821.7,
    47.3,
    58.1,
    754.3,
    1585.9,
    348.8,
    545.5,
    209.1,
    170.4,
    82.2,
    301,
    232.6,
    492.5,
    707.6,
    572.4,
    72.2,
    304.7,
    241.9,
    0.9,
    22.7
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
